$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.076.74"

$ws.Range("D3").Value = "3.815.83"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "699.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").Value = "3.815.58"
$ws.Range("E7").Value = "  -0.69%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("E10").Value = "  -0.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.53"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.90%  "

$ws.Range("E12").Value = "  +0.65%  "

$ws.Range("E13").Value = "  -1.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.82%  "

$ws.Range("D15").Value = "4.459.76"

$ws.Range("D16").Value = "3.828.91"
$ws.Range("E16").Value = "  -0.83%  "

$ws.Range("D17").Value = "71.243.30"
$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.36%  "

$ws.Range("E20").Value = "  -0.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "510.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.71"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.719"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000144"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.83%  "

$ws.Range("D26").Value = "3.969.06"
$ws.Range("E26").Value = "  -0.75%  "

$ws.Range("E27").Value = "  -0.81%  "

$ws.Range("E28").Value = "  -0.94%  "

$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.53%  "

$ws.Range("E31").Value = "  -4.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.33%  "

$ws.Range("E33").Value = "  -1.20%  "

$ws.Range("E34").Value = "  -1.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.174"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.36%  "

$ws.Range("D37").Value = "3.778.42"
$ws.Range("E37").Value = "  -0.60%  "

$ws.Range("E39").Value = "  -1.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.36%  "

$ws.Range("E42").Value = "  -1.61%  "

$ws.Range("E43").Value = "  -1.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "171.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.95%  "

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("E47").Value = "  -0.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "427.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.01%  "

$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.18%  "
